$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows whose Target cluster (column D) is "Resolving-Mac".
# These are rows 7, 13 and 19 in the original sheet; delete bottom-up so the
# row numbers of not-yet-processed rows stay stable while deleting.
$ws.Rows.Item(19).Delete()
$ws.Rows.Item(13).Delete()
$ws.Rows.Item(7).Delete()

# Refresh the NATMI-computed metric columns (G,H,I,J,M,N,O,P,Q,R,S,T) for every
# remaining data row (2-16) with the values produced by the re-run (new TPM) script.
# Values are passed as quoted numeric strings so PowerShell parses the
# scientific-notation ones (e.g. "3.54E-05") without choking on the bare "E-05" token;
# Excel COM coerces the string to a Double on assignment either way.
$ws.Range("G2").Value = "0.04110066666666667"
$ws.Range("H2").Value = "0.123302"
$ws.Range("I2").Value = "0.02671259512010182"
$ws.Range("J2").Value = "0.02671259512010182"
$ws.Range("M2").Value = "19.412944"
$ws.Range("N2").Value = "58.238832"
$ws.Range("O2").Value = "0.3138213864953257"
$ws.Range("P2").Value = "0.3138213864953256"
$ws.Range("Q2").Value = "0.7978849403626667"
$ws.Range("R2").Value = "7.180964463263999"
$ws.Range("S2").Value = "0.008382983637478626"
$ws.Range("T2").Value = "0.008382983637478622"

$ws.Range("G3").Value = "0.04110066666666667"
$ws.Range("H3").Value = "0.123302"
$ws.Range("I3").Value = "0.02671259512010182"
$ws.Range("J3").Value = "0.02671259512010182"
$ws.Range("M3").Value = "4.741326666666667"
$ws.Range("N3").Value = "14.22398"
$ws.Range("O3").Value = "0.07664626799317993"
$ws.Range("P3").Value = "0.07664626799317992"
$ws.Range("Q3").Value = "0.1948716868844445"
$ws.Range("R3").Value = "1.75384518196"
$ws.Range("S3").Value = "0.002047420724368635"
$ws.Range("T3").Value = "0.002047420724368634"

$ws.Range("G4").Value = "0.04110066666666667"
$ws.Range("H4").Value = "0.123302"
$ws.Range("I4").Value = "0.02671259512010182"
$ws.Range("J4").Value = "0.02671259512010182"
$ws.Range("M4").Value = "0.08199133333333333"
$ws.Range("N4").Value = "0.245974"
$ws.Range("O4").Value = "0.001325436982008864"
$ws.Range("P4").Value = "0.001325436982008864"
$ws.Range("Q4").Value = "0.003369898460888889"
$ws.Range("R4").Value = "0.030329086148"
$ws.Range("S4").Value = "3.540586145761246E-05"
$ws.Range("T4").Value = "3.540586145761246E-05"

$ws.Range("G5").Value = "0.04110066666666667"
$ws.Range("H5").Value = "0.123302"
$ws.Range("I5").Value = "0.02671259512010182"
$ws.Range("J5").Value = "0.02671259512010182"
$ws.Range("M5").Value = "37.49275133333333"
$ws.Range("N5").Value = "112.478254"
$ws.Range("O5").Value = "0.6060918532990739"
$ws.Range("P5").Value = "0.6060918532990738"
$ws.Range("Q5").Value = "1.540977074967556"
$ws.Range("R5").Value = "13.868793674708"
$ws.Range("S5").Value = "0.01619028628277031"
$ws.Range("T5").Value = "0.01619028628277031"

$ws.Range("G6").Value = "0.04110066666666667"
$ws.Range("H6").Value = "0.123302"
$ws.Range("I6").Value = "0.02671259512010182"
$ws.Range("J6").Value = "0.02671259512010182"
$ws.Range("M6").Value = "0.130837"
$ws.Range("N6").Value = "0.3925110000000001"
$ws.Range("O6").Value = "0.002115055230411674"
$ws.Range("P6").Value = "0.002115055230411674"
$ws.Range("Q6").Value = "0.005377487924666667"
$ws.Range("R6").Value = "0.04839739132200001"
$ws.Range("S6").Value = "5.649861402664073E-05"
$ws.Range("T6").Value = "5.649861402664073E-05"

$ws.Range("G7").Value = "0.9789586666666666"
$ws.Range("H7").Value = "2.936876"
$ws.Range("I7").Value = "0.6362555311831452"
$ws.Range("J7").Value = "0.636255531183145"
$ws.Range("M7").Value = "19.412944"
$ws.Range("N7").Value = "58.238832"
$ws.Range("O7").Value = "0.3138213864953257"
$ws.Range("P7").Value = "0.3138213864953256"
$ws.Range("Q7").Value = "19.00446977431466"
$ws.Range("R7").Value = "171.040227968832"
$ws.Range("S7").Value = "0.1996705929612146"
$ws.Range("T7").Value = "0.1996705929612144"

$ws.Range("G8").Value = "0.9789586666666666"
$ws.Range("H8").Value = "2.936876"
$ws.Range("I8").Value = "0.6362555311831452"
$ws.Range("J8").Value = "0.636255531183145"
$ws.Range("M8").Value = "4.741326666666667"
$ws.Range("N8").Value = "14.22398"
$ws.Range("O8").Value = "0.07664626799317993"
$ws.Range("P8").Value = "0.07664626799317992"
$ws.Range("Q8").Value = "4.641562831831111"
$ws.Range("R8").Value = "41.77406548648"
$ws.Range("S8").Value = "0.04876661195520639"
$ws.Range("T8").Value = "0.04876661195520637"

$ws.Range("G9").Value = "0.9789586666666666"
$ws.Range("H9").Value = "2.936876"
$ws.Range("I9").Value = "0.6362555311831452"
$ws.Range("J9").Value = "0.636255531183145"
$ws.Range("M9").Value = "0.08199133333333333"
$ws.Range("N9").Value = "0.245974"
$ws.Range("O9").Value = "0.001325436982008864"
$ws.Range("P9").Value = "0.001325436982008864"
$ws.Range("Q9").Value = "0.08026612635822222"
$ws.Range("R9").Value = "0.722395137224"
$ws.Range("S9").Value = "0.0008433166110378345"
$ws.Range("T9").Value = "0.0008433166110378343"

$ws.Range("G10").Value = "0.9789586666666666"
$ws.Range("H10").Value = "2.936876"
$ws.Range("I10").Value = "0.6362555311831452"
$ws.Range("J10").Value = "0.636255531183145"
$ws.Range("M10").Value = "37.49275133333333"
$ws.Range("N10").Value = "112.478254"
$ws.Range("O10").Value = "0.6060918532990739"
$ws.Range("P10").Value = "0.6060918532990738"
$ws.Range("Q10").Value = "36.70385385494489"
$ws.Range("R10").Value = "330.3346846945039"
$ws.Range("S10").Value = "0.3856292940665791"
$ws.Range("T10").Value = "0.385629294066579"

$ws.Range("G11").Value = "0.9789586666666666"
$ws.Range("H11").Value = "2.936876"
$ws.Range("I11").Value = "0.6362555311831452"
$ws.Range("J11").Value = "0.636255531183145"
$ws.Range("M11").Value = "0.130837"
$ws.Range("N11").Value = "0.3925110000000001"
$ws.Range("O11").Value = "0.002115055230411674"
$ws.Range("P11").Value = "0.002115055230411674"
$ws.Range("Q11").Value = "0.1280840150706667"
$ws.Range("R11").Value = "1.152756135636"
$ws.Range("S11").Value = "0.001345715589107269"
$ws.Range("T11").Value = "0.001345715589107269"

$ws.Range("G12").Value = "0.5185656666666667"
$ws.Range("H12").Value = "1.555697"
$ws.Range("I12").Value = "0.3370318736967531"
$ws.Range("J12").Value = "0.3370318736967531"
$ws.Range("M12").Value = "19.412944"
$ws.Range("N12").Value = "58.238832"
$ws.Range("O12").Value = "0.3138213864953257"
$ws.Range("P12").Value = "0.3138213864953256"
$ws.Range("Q12").Value = "10.06688624732267"
$ws.Range("R12").Value = "90.601976225904"
$ws.Range("S12").Value = "0.1057678098966326"
$ws.Range("T12").Value = "0.1057678098966325"

$ws.Range("G13").Value = "0.5185656666666667"
$ws.Range("H13").Value = "1.555697"
$ws.Range("I13").Value = "0.3370318736967531"
$ws.Range("J13").Value = "0.3370318736967531"
$ws.Range("M13").Value = "4.741326666666667"
$ws.Range("N13").Value = "14.22398"
$ws.Range("O13").Value = "0.07664626799317993"
$ws.Range("P13").Value = "0.07664626799317992"
$ws.Range("Q13").Value = "2.458689223784444"
$ws.Range("R13").Value = "22.12820301406"
$ws.Range("S13").Value = "0.02583223531360491"
$ws.Range("T13").Value = "0.0258322353136049"

$ws.Range("G14").Value = "0.5185656666666667"
$ws.Range("H14").Value = "1.555697"
$ws.Range("I14").Value = "0.3370318736967531"
$ws.Range("J14").Value = "0.3370318736967531"
$ws.Range("M14").Value = "0.08199133333333333"
$ws.Range("N14").Value = "0.245974"
$ws.Range("O14").Value = "0.001325436982008864"
$ws.Range("P14").Value = "0.001325436982008864"
$ws.Range("Q14").Value = "0.04251789043088889"
$ws.Range("R14").Value = "0.382661013878"
$ws.Range("S14").Value = "0.000446714509513417"
$ws.Range("T14").Value = "0.0004467145095134169"

$ws.Range("G15").Value = "0.5185656666666667"
$ws.Range("H15").Value = "1.555697"
$ws.Range("I15").Value = "0.3370318736967531"
$ws.Range("J15").Value = "0.3370318736967531"
$ws.Range("M15").Value = "37.49275133333333"
$ws.Range("N15").Value = "112.478254"
$ws.Range("O15").Value = "0.6060918532990739"
$ws.Range("P15").Value = "0.6060918532990738"
$ws.Range("Q15").Value = "19.44245359033756"
$ws.Range("R15").Value = "174.982082313038"
$ws.Range("S15").Value = "0.2042722729497245"
$ws.Range("T15").Value = "0.2042722729497244"

$ws.Range("G16").Value = "0.5185656666666667"
$ws.Range("H16").Value = "1.555697"
$ws.Range("I16").Value = "0.3370318736967531"
$ws.Range("J16").Value = "0.3370318736967531"
$ws.Range("M16").Value = "0.130837"
$ws.Range("N16").Value = "0.3925110000000001"
$ws.Range("O16").Value = "0.002115055230411674"
$ws.Range("P16").Value = "0.002115055230411674"
$ws.Range("Q16").Value = "0.06784757612966667"
$ws.Range("R16").Value = "0.6106281851670001"
$ws.Range("S16").Value = "0.0007128410272777644"
$ws.Range("T16").Value = "0.0007128410272777644"
